$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("D1").Value = "ObjectName"

# New column values (row 2 and row 3 left blank, as in the diff)
$ws.Range("D4").Value  = "default"
$ws.Range("D5").Value  = "default"
$ws.Range("D6").Value  = "default"
$ws.Range("D7").Value  = "default"
$ws.Range("D8").Value  = "player"
$ws.Range("D9").Value  = "default"
$ws.Range("D10").Value = "default"
$ws.Range("D11").Value = "default"
$ws.Range("D12").Value = "default"
$ws.Range("D13").Value = "default"

# Update selection to match target (L13)
$ws.Range("L13").Select()
